$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A (stored width goes from 16.42578125 to 15.42578125; the COM
# model snaps ColumnWidth to 1/6-character increments, so 15.5 is the closest
# representable width - one character narrower, matching the diff's intent).
$ws.Columns.Item(1).ColumnWidth = 14.6

$values = @(
    0.12818632190988666,
    -0.0059999999887025979,
    -0.0039999999908015837,
    -0.0079999999827400359,
    -0.0029999999910765851,
    -0.0019999999910229604,
    -0.0099999999761464231,
    -0.0099999999755810975,
    -0.0019999999896285203,
    -0.001999999989154233,
    0.0032038245706713653,
    -0.0034999999862481701,
    -0.0034999999857125985,
    -0.0079999999773514574,
    0.0097805020665635567,
    -0.001999999988138601,
    -0.001999999987854828,
    -0.0039999999841100475,
    -0.0039999999923572283,
    -0.0039999999917625928,
    -0.0039999999916924267,
    -0.003999999991604497,
    0.012062290698499289,
    -0.019999999960004899,
    -0.019999999959466663,
    -0.00249999999000039,
    -0.0024999999898067671,
    -0.0019999999900308651,
    -0.0069999999804570834,
    -0.059999999883465271,
    -0.0069999999821046544,
    0.063732523984885603,
    -0.0039999999880784287
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
